$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A128").Value = "idRestorePurcases"
$ws.Range("B128").Value = "You may restore your purchases by tapping this button."

$ws.Range("A129").Value = "idBackToPuzzle"
$ws.Range("B129").Value = "This button returns you to the puzzle."

$ws.Range("C129").Value = "Эта кнопка возвращает\nк головоломке."
$ws.Range("C128").Value = "Эта кнопка восстанавливает\nВаши покупки."

$ws.Range("A124:C124").Copy() | Out-Null
$ws.Range("A128:C128").PasteSpecial(-4122) | Out-Null

$ws.Range("A123:C123").Copy() | Out-Null
$ws.Range("A129:C129").PasteSpecial(-4122) | Out-Null

$ws.Range("A129").Select() | Out-Null
